$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.100.30'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.261.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '546.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.98%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -0.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.41'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.54%  '
$ws.Range("E10").Value = '  +3.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.433'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.825.71'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.29%  '
$ws.Range("E13").Value = '  -1.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.54%  '
$ws.Range("E15").Value = '  +2.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.119.88'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.258.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.42%  '
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("E19").Value = '  +3.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '379.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.534'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("E25").Value = '  +1.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.66'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0931'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.94%  '
$ws.Range("E29").Value = '  +2.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.67'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.95%  '
$ws.Range("E31").Value = '  +1.86%  '
$ws.Range("E32").Value = '  +2.79%  '
$ws.Range("E33").Value = '  +7.34%  '
$ws.Range("E34").Value = '  +4.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.70'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.45'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.814.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0723'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.25%  '
$ws.Range("E40").Value = '  +1.08%  '
$ws.Range("E41").Value = '  +6.88%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.734'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.306.52'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.41%  '
$ws.Range("E46").Value = '  +2.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.02'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.61%  '
$ws.Range("E49").Value = '  +0.58%  '
$ws.Range("E50").Value = '  +5.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '278.89'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.42%  '
